$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_rerun_results_tracker")
$ws.Range("J186").NumberFormat = "0.00"
Write-Host $ws.Range("J186").NumberFormat
Write-Host $ws.Range("J186").Value
